# edit.ps1 - reproduce the target commit:
#   1) Re-style the three tables (slides 14-16) from table style
#      {219DEB86-0FA9-41FD-983D-CBA69F6D9310} to
#      {7848CB47-DF85-444D-BA4C-B9AD63EB1DFD}.
#   2) Swap the presentation's colour theme from the "Integral" (Red
#      Violet) palette to the default Office palette (the values that
#      used to live in theme2.xml / the Notes Master slot).

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$oldStyle = "{219DEB86-0FA9-41FD-983D-CBA69F6D9310}"
$newStyle = "{7848CB47-DF85-444D-BA4C-B9AD63EB1DFD}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyle)
        }
    }
}

# --- 2) Theme colours --------------------------------------------------
# Helper: turn an RRGGBB hex string into the BGR-packed decimal that the
# ColorFormat.RGB / ThemeColor.RGB properties expect (classic
# Windows COLORREF / VBA RGB() ordering).
function Convert-HexToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office default theme palette, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($idx = 1; $idx -le $officeThemeColors.Count; $idx++) {
    $themeColorScheme.Colors($idx).RGB = Convert-HexToBgr $officeThemeColors[$idx - 1]
}
